$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 252, shifting existing rows 252:328 down to 253:329
$ws.Rows.Item(252).Insert()

# Populate the newly inserted row 252 with the new record's data
$ws.Range("A252").Value = 8
$ws.Range("B252").Value = "Terminal La Palmera de La Serena"
$ws.Range("C252").Value = "Coquimbo"
$ws.Range("D252").Value = 44964
$ws.Range("E252").Value = 4
$ws.Range("F252").Value = "Fruta"
$ws.Range("G252").Value = 100103
$ws.Range("H252").Value = "Frutos de hueso (carozo)"
$ws.Range("I252").Value = 100103002
$ws.Range("J252").Value = "Ciruela"
$ws.Range("K252").Value = "Larry Ann"
$ws.Range("L252").Value = "Primera"
$ws.Range("M252").Value = 12
$ws.Range("N252").Value = 315000
$ws.Range("O252").Value = 320000
$ws.Range("P252").Value = 317500
$ws.Range("Q252").Value = "$/bins (450 kilos)"
$ws.Range("R252").Value = "Región de O'Higgins"
$ws.Range("S252").Value = 706
$ws.Range("T252").Value = 450
